$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counts: Cant. Trabajadores / Cant. Periodos go from 2 to 1
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Update VALOR MORA summary figure
$ws.Range("E11").Value = 29509

# Replace the first worker's data (row 16) with the second worker's data
# (previously in row 17), since the first worker is being removed.
$ws.Range("C16").Value = "1143373645"
$ws.Range("D16").Value = "ALEJANDRO JOSE LEON BLANCO"
$ws.Range("E16").Value = "1711"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 781242

# Remove the now-duplicated second worker's row (row 17), shifting the trailing
# signature rows up by one (old rows 22/23 become 21/22).
$ws.Rows("17").Delete()
